# Add two new columns (A and B) at the start of the sheet, shifting the
# existing data (originally in columns A:O) two columns to the right
# (becoming C:Q). Populate the new columns with a "col"/"tar" header and
# "STAR"/"pp" values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column A; this shifts all existing
# columns (and their formulas) two positions to the right automatically.
$ws.Range("A:B").Insert()

# Header row for the two new columns.
$ws.Range("A1").Value = "col"
$ws.Range("B1").Value = "tar"

# Data rows (originally rows 2-7) get the same label in each new column.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "STAR"
    $ws.Cells.Item($r, 2).Value = "pp"
}

# Match the centered-horizontal-alignment style used by the rest of the
# data table (same style as the columns that got shifted right).
$ws.Range("A1:B7").HorizontalAlignment = -4108

# Re-apply the two formula columns (now P and Q, previously N and O) as a
# single range assignment so the engine keeps writing them out as shared
# formulas (t="shared") across rows 3-7, matching the original layout.
$ws.Range("P3:P7").Formula = "=SQRT(K3*K3-N3*N3-O3*O3)"
$ws.Range("Q3:Q7").Formula = "=SQRT(N3*N3+O3*O3)"

# Update the active selection to match the post-edit state.
$ws.Range("E15").Select()
